$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 9.380719000000001
$ws.Range("H2").Value = 28.142157
$ws.Range("I2").Value = 0.03679977590837273
$ws.Range("J2").Value = 0.03679977590837273
$ws.Range("Q2").Value = 0.5853099620050001
$ws.Range("R2").Value = 5.267789658045
$ws.Range("S2").Value = 0.03679977590837273
$ws.Range("T2").Value = 0.03679977590837273

# Row 3
$ws.Range("I3").Value = 0.3547860986448385
$ws.Range("J3").Value = 0.3547860986448385
$ws.Range("S3").Value = 0.3547860986448385
$ws.Range("T3").Value = 0.3547860986448385

# Row 4
$ws.Range("G4").Value = 100.179423
$ws.Range("H4").Value = 300.538269
$ws.Range("I4").Value = 0.3929954960840508
$ws.Range("J4").Value = 0.3929954960840508
$ws.Range("Q4").Value = 6.250695098085
$ws.Range("R4").Value = 56.256255882765
$ws.Range("S4").Value = 0.3929954960840508
$ws.Range("T4").Value = 0.3929954960840508

# Row 5
$ws.Range("G5").Value = 1.427630666666667
$ws.Range("H5").Value = 4.282892
$ws.Range("I5").Value = 0.005600475679236752
$ws.Range("J5").Value = 0.005600475679236752
$ws.Range("Q5").Value = 0.08907701544666667
$ws.Range("R5").Value = 0.8016931390200001
$ws.Range("S5").Value = 0.005600475679236752
$ws.Range("T5").Value = 0.005600475679236752

# Row 6
$ws.Range("G6").Value = 53.48524799999999
$ws.Range("H6").Value = 160.455744
$ws.Range("I6").Value = 0.2098181536835013
$ws.Range("J6").Value = 0.2098181536835013
$ws.Range("Q6").Value = 3.337212048959999
$ws.Range("R6").Value = 30.03490844064
$ws.Range("S6").Value = 0.2098181536835013
$ws.Range("T6").Value = 0.2098181536835013
